# Apply "data updated till 12Jan 8AM" edits to the Retailer-wise order sheet.
# Only the underlying daily (H:AL) entries need to change by hand - the
# E/F columns (per-row totals) and the whole row-2 grand-total row are
# formulas and recompute automatically once the source cells below change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> column -> new value for the newly-entered daily figures
$ws.Range("R4").Value  = 1560
$ws.Range("R5").Value  = 1040
$ws.Range("S6").Value  = 5200
$ws.Range("S12").Value = 1040
$ws.Range("R23").Value = 3120
$ws.Range("S25").Value = 3120
$ws.Range("S29").Value = 5200
$ws.Range("R40").Value = 2080
$ws.Range("R41").Value = 3120
$ws.Range("R48").Value = 3120
$ws.Range("R52").Value = 1040
$ws.Range("R54").Value = 2080
$ws.Range("R70").Value = 3120
$ws.Range("S71").Value = 2080
$ws.Range("R72").Value = 2080
$ws.Range("S82").Value = 5200

# R4 picks up the "new entry" highlight fill used elsewhere in the sheet
# (solid fill rgb FFE1B884), matching the other freshly filled-in cells.
$ws.Range("R4").Interior.Color = 8698081

# Reflect where the author had scrolled/selected before saving.
$ws.Range("O45").Select()
